# Insert a new weekly price record for "Perejil" (Vega Modelo de Temuco)
# right before the current row 160, shifting every following row down by
# one (row 259 -> 260), and extending the used range to A1:R260.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 160; everything from 160..259 shifts to 161..260.
$ws.Rows("160").Insert(-4121)   # -4121 = xlShiftDown

# Populate the newly inserted row with the new record's data.
$ws.Range("A160").Value = 10
$ws.Range("B160").Value = 'Vega Modelo de Temuco'
$ws.Range("C160").Value = 'La Araucanía'
$ws.Range("D160").Value = 44582
$ws.Range("E160").Value = 9
$ws.Range("F160").Value = 100112044
$ws.Range("G160").Value = 'Perejil'
$ws.Range("H160").Value = 'Sin especificar'
$ws.Range("I160").Value = 'Primera'
$ws.Range("J160").Value = 40
$ws.Range("K160").Value = 5000
$ws.Range("L160").Value = 5000
$ws.Range("M160").Value = 5000
$ws.Range("N160").Value = '$/docena de atados (3 kilos)'
$ws.Range("O160").Value = 'Provincia de Cautín'
$ws.Range("P160").Value = 1667
$ws.Range("Q160").Value = 3
$ws.Range("R160").Value = 'Hortaliza'
